# Add a "Total:" row beneath the product backlog data (rows 2-12),
# summing each numeric column (B-F) with a SUM formula, matching the
# bold header style already used in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label
$ws.Range("A13").Value = "Total:"

# Column totals (individual formulas, one per cell, so each carries its
# own <f> element rather than a shared-formula group)
$ws.Range("B13").Formula = "=SUM(B2:B12)"
$ws.Range("C13").Formula = "=SUM(C2:C12)"
$ws.Range("D13").Formula = "=SUM(D2:D12)"
$ws.Range("E13").Formula = "=SUM(E2:E12)"
$ws.Range("F13").Formula = "=SUM(F2:F12)"

# Match the bold styling used by the header row
$ws.Range("A13:F13").Font.Bold = $true

# Move the active selection, as happened in the authored edit
$ws.Range("C15").Select() | Out-Null
